# Updated cryptos list on Wed Feb 22 15:50:14 UTC 2023 with GitHub Actions
# Applies the scraped price/volume refresh + two coin-row reorders (36<->37, 41<->42).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking text values are written with a leading apostrophe so Excel
# keeps them as text (matching the source sheet, which stores Price/Volume as
# inline strings, not numbers).

# Row 2
$ws.Range("D2").Value = "23.897.91"
$ws.Range("E2").Value = "  -3.10%  "

# Row 3
$ws.Range("D3").Value = "1.625.18"
$ws.Range("E3").Value = "  -3.07%  "

# Row 4
$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  +0.43%  "

# Row 5
$ws.Range("D5").Value = "'1.006"
$ws.Range("E5").Value = "  +0.37%  "

# Row 6
$ws.Range("D6").Value = "'306.75"
$ws.Range("E6").Value = "  -2.31%  "

# Row 7
$ws.Range("D7").Value = "'0.3907"
$ws.Range("E7").Value = "  +0.36%  "

# Row 8
$ws.Range("D8").Value = "'0.3814"
$ws.Range("E8").Value = "  -3.19%  "

# Row 9
$ws.Range("D9").Value = "'1.006"
$ws.Range("E9").Value = "  +0.38%  "

# Row 10
$ws.Range("D10").Value = "'49.66"
$ws.Range("E10").Value = "  -4.35%  "

# Row 11
$ws.Range("D11").Value = "'1.359"
$ws.Range("E11").Value = "  -2.38%  "

# Row 12
$ws.Range("D12").Value = "'0.08431"
$ws.Range("E12").Value = "  -2.40%  "

# Row 13
$ws.Range("D13").Value = "'23.71"
$ws.Range("E13").Value = "  -5.81%  "

# Row 14
$ws.Range("D14").Value = "'6.985"
$ws.Range("E14").Value = "  -4.39%  "

# Row 15
$ws.Range("D15").Value = "'0.00001269"
$ws.Range("E15").Value = "  -3.35%  "

# Row 16
$ws.Range("D16").Value = "'7.418"
$ws.Range("E16").Value = "  -4.57%  "

# Row 17
$ws.Range("D17").Value = "1.629.09"
$ws.Range("E17").Value = "  -3.58%  "

# Row 18
$ws.Range("D18").Value = "'92.68"
$ws.Range("E18").Value = "  -1.05%  "

# Row 19
$ws.Range("D19").Value = "'0.06909"
$ws.Range("E19").Value = "  -2.11%  "

# Row 20
$ws.Range("D20").Value = "'19.80"
$ws.Range("E20").Value = "  -3.97%  "

# Row 21
$ws.Range("D21").Value = "'6.834"
$ws.Range("E21").Value = "  -3.20%  "

# Row 22
$ws.Range("E22").Value = "  +0.25%  "

# Row 23
$ws.Range("D23").Value = "'13.35"
$ws.Range("E23").Value = "  -4.60%  "

# Row 24
$ws.Range("D24").Value = "23.909.97"
$ws.Range("E24").Value = "  -3.07%  "

# Row 25
$ws.Range("D25").Value = "'2.393"
$ws.Range("E25").Value = "  +1.39%  "

# Row 26
$ws.Range("D26").Value = "'2.872"
$ws.Range("E26").Value = "  +5.34%  "

# Row 27
$ws.Range("D27").Value = "'22.09"
$ws.Range("E27").Value = "  -4.52%  "

# Row 28
$ws.Range("D28").Value = "'157.67"
$ws.Range("E28").Value = "  -2.87%  "

# Row 29
$ws.Range("D29").Value = "'138.62"
$ws.Range("E29").Value = "  -5.39%  "

# Row 30
$ws.Range("D30").Value = "'5.244"
$ws.Range("E30").Value = "  -8.71%  "

# Row 31
$ws.Range("D31").Value = "'7.645"
$ws.Range("E31").Value = "  -3.16%  "

# Row 32
$ws.Range("D32").Value = "'2.464"
$ws.Range("E32").Value = "  -2.75%  "

# Row 33
$ws.Range("D33").Value = "1.802.02"
$ws.Range("E33").Value = "  -7.00%  "

# Row 34
$ws.Range("D34").Value = "'0.07925"
$ws.Range("E34").Value = "  -5.31%  "

# Row 35
$ws.Range("D35").Value = "'0.02877"
$ws.Range("E35").Value = "  -5.27%  "

# Row 36
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'0.9582"
$ws.Range("E36").Value = "  -2.49%  "

# Row 37
$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").Value = "'6.582"
$ws.Range("E37").Value = "  -4.00%  "

# Row 38
$ws.Range("D38").Value = "'0.2652"
$ws.Range("E38").Value = "  -5.88%  "

# Row 39
$ws.Range("D39").Value = "'0.09131"
$ws.Range("E39").Value = "  -3.72%  "

# Row 40
$ws.Range("D40").Value = "'10.21"
$ws.Range("E40").Value = "  -3.08%  "

# Row 41
$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").Value = "'13.40"
$ws.Range("E41").Value = "  -0.96%  "

# Row 42
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'1.420"
$ws.Range("E42").Value = "  -8.41%  "

# Row 43
$ws.Range("D43").Value = "'0.7423"
$ws.Range("E43").Value = "  -6.04%  "

# Row 44
$ws.Range("D44").Value = "'15.91"
$ws.Range("E44").Value = "  -3.95%  "

# Row 45
$ws.Range("D45").Value = "'0.6817"
$ws.Range("E45").Value = "  -4.31%  "

# Row 46
$ws.Range("D46").Value = "'2.434"
$ws.Range("E46").Value = "  -4.95%  "

# Row 47
$ws.Range("D47").Value = "'4.064"
$ws.Range("E47").Value = "  -2.94%  "

# Row 48
$ws.Range("D48").Value = "'1.005"
$ws.Range("E48").Value = "  +0.30%  "

# Row 49
$ws.Range("D49").Value = "'0.08245"
$ws.Range("E49").Value = "  -4.72%  "

# Row 50
$ws.Range("D50").Value = "'132.49"
$ws.Range("E50").Value = "  -3.45%  "

# Row 51
$ws.Range("D51").Value = "'1.249"
$ws.Range("E51").Value = "  -5.58%  "
